$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 555
$ws.Range("F3").Value = 248
$ws.Range("F4").Value = 580
$ws.Range("F5").Value = 1372
$ws.Range("F6").Value = 682
$ws.Range("F7").Value = 362
$ws.Range("F8").Value = 40
$ws.Range("F9").Value = 160
$ws.Range("F11").Value = 6392
$ws.Range("F12").Value = 120
$ws.Range("F13").Value = 28
$ws.Range("F15").Value = 4754
$ws.Range("F19").Value = 5560
$ws.Range("F20").Value = 7400
$ws.Range("F22").Value = 1094
$ws.Range("F23").Value = 767
$ws.Range("F24").Value = 4049
$ws.Range("F25").Value = 568
$ws.Range("F27").Value = 235
$ws.Range("F29").Value = 1068
$ws.Range("F30").Value = 1502
$ws.Range("F31").Value = 573
$ws.Range("F32").Value = 701
$ws.Range("F35").Value = 1930
$ws.Range("F36").Value = 239
$ws.Range("F38").Value = 1257
$ws.Range("F40").Value = 702
$ws.Range("F41").Value = 327
$ws.Range("F42").Value = 1724
$ws.Range("F43").Value = 3721
$ws.Range("F45").Value = 351
$ws.Range("F46").Value = 451
$ws.Range("F47").Value = 27
$ws.Range("F48").Value = 103
$ws.Range("F49").Value = 3967
# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 1283
$ws.Range("F5").Value = 46
$ws.Range("F29").Value = 87
# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 4496
# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 4496
$ws.Range("F4").Value = 555
$ws.Range("F5").Value = 1283
$ws.Range("F8").Value = 248
$ws.Range("F9").Value = 580
$ws.Range("F11").Value = 1372
$ws.Range("F13").Value = 682
$ws.Range("F14").Value = 362
$ws.Range("F15").Value = 40
$ws.Range("F16").Value = 160
$ws.Range("C18").Value = "杭州·AD04动漫展"
$ws.Range("E18").Value = "2024.07.13 10:00-07.14 17:00"
$ws.Range("F18").Value = 6392
$ws.Range("G18").Value = 75
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=85012"
$ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202405/y1iKqqnh1715326769523.jpeg"
$ws.Range("F20").Value = 4754
$ws.Range("F21").Value = 5560
$ws.Range("F22").Value = 5560
$ws.Range("C23").Value = "杭州·TCD国潮动漫游戏嘉年华"
$ws.Range("D23").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E23").Value = "2024.07.20 09:30-07.21 17:00"
$ws.Range("F23").Value = 7400
$ws.Range("G23").Value = 65
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=85699"
$ws.Range("I23").Value = "//i1.hdslb.com/bfs/openplatform/202406/QzaksReK1718190369702.jpeg"
$ws.Range("C24").Value = "杭州·次元幻想--二次元全女夜场"
$ws.Range("D24").Value = "保淑路2号 The Queen皇后"
$ws.Range("E24").Value = "2024.07.20 13:00-07.20 19:00"
$ws.Range("F24").Value = 1094
$ws.Range("G24").Value = 158
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=81808"
$ws.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202406/KEJ8V0WB1718178102884.jpeg"
$ws.Range("C25").Value = "杭州·生如夏花国乙only·日夜场"
$ws.Range("D25").Value = "祥符街道花园岗街181号 格乐利雅婚礼艺术中心(天空之城店)"
$ws.Range("E25").Value = "2024.07.20 10:00-07.20 22:30"
$ws.Range("F25").Value = 767
$ws.Range("G25").Value = 135
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=85496"
$ws.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202405/Qut2ZdAi1715411977772.jpeg"
$ws.Range("C26").Value = "杭州·第五届华盟次元嘉年华&周年庆狂欢"
$ws.Range("D26").Value = "创意路1号 中国智谷富春园区"
$ws.Range("E26").Value = "2024.07.20 10:00-07.21 17:00"
$ws.Range("F26").Value = 4049
$ws.Range("G26").Value = 58
$ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=84762"
$ws.Range("I26").Value = "//i0.hdslb.com/bfs/openplatform/202404/uE6OVg6T1713885553204.jpeg"
$ws.Range("C27").Value = "杭州·第四届ArknightsOnly·狼与黑荆棘（明日方舟Only）"
$ws.Range("D27").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E27").Value = "2024.07.20 10:00-07.20 17:00"
$ws.Range("F27").Value = 568
$ws.Range("G27").Value = 79
$ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=86305"
$ws.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202405/cpoiCink1716554216810.png"
$ws.Range("F30").Value = 1068
$ws.Range("F31").Value = 1502
$ws.Range("F32").Value = 573
$ws.Range("F33").Value = 701
$ws.Range("F35").Value = 242
$ws.Range("F36").Value = 1930
$ws.Range("F41").Value = 702
$ws.Range("F42").Value = 327
$ws.Range("F43").Value = 87
$ws.Range("F44").Value = 3721
$ws.Range("F47").Value = 351
$ws.Range("F48").Value = 103
$ws.Range("F50").Value = 3967

Write-Output "done"